# Generate Report for Handback
# Updates the localization-status report: flips the Overview "Status" column
# from "Ready for handoff" to "Handed back: in sync with en-US", refreshes the
# "Latest Handback DateTime" stamps for zh-cn/de-de now that the files came
# back in sync, and clears the (now stale) "Error Detail" messages.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!Status (zh-cn / de-de columns) for both rows
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn!Latest Handback DateTime refreshed, Error Detail cleared
$wsZhCn.Range("K2").Value = "2016-11-15 16:44:57"
$wsZhCn.Range("K3").Value = "2016-11-15 16:44:57"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

# de-de!Latest Handback DateTime refreshed, Error Detail cleared
$wsDeDe.Range("K2").Value = "2016-11-15 16:45:17"
$wsDeDe.Range("K3").Value = "2016-11-15 16:45:17"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

# Column widths widen to fit the longer "Handed back" status text, and the
# Error Detail column narrows now that it's empty.
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

$wsZhCn.Columns.Item(3).ColumnWidth = 29.15
$wsZhCn.Columns.Item(16).ColumnWidth = 12.85

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15
$wsDeDe.Columns.Item(16).ColumnWidth = 12.85
